$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before I to store the initial water level aggregation type.
# This used to default to "max" implicitly; now it is read from sqlite and
# defaults to "min", so make sure we persist a flag for it here.
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "water_level_ini_type"

# Flag the rows that already define an initial water level raster (rows 4-7)
$ws.Range("I4:I7").Value = 1

# Restore the selection state of the sheet
$ws.Range("Z5").Select() | Out-Null
